$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.2762
$ws.Range("D3").Value = -7.249599999999997
$ws.Range("E5").Value = 12.61159999999999
$ws.Range("D14").Value = -8.088300000000002
$ws.Range("D16").Value = -8.195299999999996
$ws.Range("E16").Value = 12.97230000000001
$ws.Range("D21").Value = -7.524999999999995
$ws.Range("D23").Value = -7.549099999999997
$ws.Range("D25").Value = -8.117999999999999
